# Docs/Tasks.xlsx - "Added admin menu page, Users contoller"
#
# Target sheet: "M0 - Account Mgmt" (3rd worksheet)
#
# The edit inserts two new task rows into the "ADMIN FEATURES" block
# (for the new Admin Game controller + user-role listing work), which
# pushes the rest of the sheet (USER FEATURES / MISC blocks) down by
# two rows, and then tweaks a handful of individual cells' text/status.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("M0 - Account Mgmt")

# --- Make room for 2 new rows right after the existing "ADMIN FEATURES"
#     rows (old rows 7 & 8 -> new rows 9 & 10). Excel-style row insert
#     shifts everything below down and keeps the moved rows' formatting
#     (e.g. the "USER FEATURES"/"MISC" section headers, the Input-style
#     status cells, etc.) intact automatically.
$ws.Range("A7:A8").EntireRow.Insert()

# --- Row 2: Create Admin View -> now Done, with a note about how it
#     was implemented (restricted by role).
$ws.Range("C2").Value = "Done"
$ws.Range("C2").Style = "Good"
$ws.Range("D2").Value = "Admin view created and restricted by role. "

# --- Row 4: Change User password actually belongs to the User
#     Controller (not Admin Controller); status note updated to
#     reflect the real current limitation.
$ws.Range("B4").Value = "User Controller"
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = "Can edit the hash, but don't have a way to manually reset the password. "
$ws.Range("D4").Style = "Neutral"

# --- Row 5: Delete User account also belongs to User Controller and
#     is now Done.
$ws.Range("B5").Value = "User Controller"
$ws.Range("C5").Value = "Done"
$ws.Range("C5").Style = "Good"

# --- Row 6: Promote User to Admin note on what's needed next.
$ws.Range("D6").Value = "Create additional viewmodel first. "

# --- New row 7: Create Admin Game controller (Game Controller), done.
$ws.Range("A7").Value = "Create Admin Game controller"
$ws.Range("B7").Value = "Game Controller"
$ws.Range("C7").Value = "Done"
$ws.Range("C7").Style = "Good"

# --- New row 8: Modify Users list to display roles, still NEXT, needs
#     a new ViewModel.
$ws.Range("A8").Value = "Modify Users list to display their roles (admin, TO, etc.) "
$ws.Range("B8").Value = "User Controller"
$ws.Range("C8").Value = "NEXT"
$ws.Range("D8").Value = "Requires a new ViewModel. "

# --- Row 11 (old "Register" row, shifted down): drop the stray note
#     that was accidentally duplicated here.
$ws.Range("D11").ClearContents()

# Restore the worksheet's selected cell (A6) as recorded in the saved
# file.
$ws.Range("A6").Select()
